$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = "4 octets"
$ws.Range("F8").Value = 80
$ws.Range("H8").Value = "Fréquence de la vibration (dans ce projet, vitesse de clignotement de la LED1)"
$ws.Range("H10").Value = "Chaîne de caracteres utilisée pour debug"
$ws.Range("F12").Value = '"Key Finder 42"'

# Normalize the leftover "no-alignment" style used by the empty bordered
# cells in rows 3-5 so it collapses onto the same style as the bordered
# header-like cells (matches a real Excel re-save, which drops the
# redundant/default <alignment> info and merges the two identical xfs).
$ws.Range("B3:H5").WrapText = $false
